# Characters.xlsx - "Add files via upload" edit
# Updates the per-character stat table on the single worksheet and moves
# the saved selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing text storage for
# numeric-looking strings (so "187" is kept as text, not coerced to the
# number 187 the way a bare Range.Value assignment would do).
# ---------------------------------------------------------------------
function Set-TextValue {
    param($addr, $text)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
}

function Set-NumValue {
    param($addr, $num)
    $ws.Range($addr).Value = $num
}

# ---------------------------------------------------------------------
# Row 4 - HP
# ---------------------------------------------------------------------
Set-TextValue "D4"  "187"
Set-TextValue "E4"  "187"
Set-TextValue "H4"  "179"
Set-TextValue "I4"  "179"
Set-TextValue "L4"  "164"
Set-TextValue "M4"  "164"
Set-TextValue "AB4" "200"
Set-TextValue "AC4" "200"

# ---------------------------------------------------------------------
# Row 5 - MP
# ---------------------------------------------------------------------
Set-TextValue "D5"  "26"
Set-TextValue "E5"  "26"
Set-TextValue "H5"  "30"
Set-TextValue "I5"  "30"
Set-TextValue "L5"  "30"
Set-TextValue "M5"  "30"
Set-TextValue "AB5" "32"
Set-TextValue "AC5" "32"

# ---------------------------------------------------------------------
# Row 6 - Atk (numeric cells)
# ---------------------------------------------------------------------
Set-NumValue "D6"  24
Set-NumValue "H6"  10
Set-NumValue "L6"  8
Set-NumValue "AB6" 90

# ---------------------------------------------------------------------
# Row 8 - Level
# ---------------------------------------------------------------------
Set-TextValue "D8"  "10"
Set-TextValue "H8"  "10"
Set-TextValue "L8"  "10"
Set-TextValue "AB8" "10"

# ---------------------------------------------------------------------
# Row 9 - PWR
# ---------------------------------------------------------------------
Set-TextValue "D9"  "17"
Set-TextValue "H9"  "4"
Set-TextValue "L9"  "4"
Set-TextValue "AB9" "18"

# ---------------------------------------------------------------------
# Row 11 - HIT
# ---------------------------------------------------------------------
Set-TextValue "D11"  "10"
Set-TextValue "H11"  "11"
Set-TextValue "L11"  "10"
Set-TextValue "AB11" "19"

# ---------------------------------------------------------------------
# Row 12 - EV
# ---------------------------------------------------------------------
Set-TextValue "D12"  "11"
Set-TextValue "H12"  "9"
Set-TextValue "L12"  "8"
Set-TextValue "AB12" "15"

# ---------------------------------------------------------------------
# Row 13 - STM
# ---------------------------------------------------------------------
Set-TextValue "D13"  "22"
Set-TextValue "H13"  "13"
Set-TextValue "L13"  "11"
Set-TextValue "AB13" "16"

# ---------------------------------------------------------------------
# Row 14 - MAG
# ---------------------------------------------------------------------
Set-TextValue "D14"  "8"
Set-TextValue "H14"  "17"
Set-TextValue "L14"  "15"
Set-TextValue "AB14" "20"

# ---------------------------------------------------------------------
# Row 15 - MDEF
# ---------------------------------------------------------------------
Set-TextValue "D15"  "16"
Set-TextValue "H15"  "22"
Set-TextValue "L15"  "21"
Set-TextValue "AB15" "24"

# ---------------------------------------------------------------------
# Row 16 - EXP (numeric cells, previously blank)
# ---------------------------------------------------------------------
Set-NumValue "D16" 3
Set-NumValue "H16" 3
Set-NumValue "L16" 3
Set-NumValue "P16" 0
Set-NumValue "T16" 0
Set-NumValue "X16" 0

# ---------------------------------------------------------------------
# Move the saved selection to AB16 (also scroll so column N is leftmost,
# to the extent the host window object supports it).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 14
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AB16").Select() | Out-Null
